# Thunderstone Capital hourLog.xlsx update
# Adds a new logged-hours entry and splits it into three rows, and
# renames the "Website creation" work description to "Building the site".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells in place -----------------------------------
# Row 6's work description was renamed from "Website creation" to
# "Building the site".
$ws.Range("D6").Value = "Building the site"

# Row 7's hour log entry was changed from "10:05-00:00" to "10:05-1:30, ".
$ws.Range("B7").Value = "10:05-1:30, "

# --- Finish filling out row 7 (total hours + work description) --------
$ws.Range("C2").Copy($ws.Range("C7"))
$ws.Range("C7").Value = 0.14583333333333334
$ws.Range("D7").Value = "Building the site"

# --- Add row 8: Nov 7, 2019 --------------------------------------------
$ws.Range("A2").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 43776
$ws.Range("B7").Copy($ws.Range("B8"))
$ws.Range("B8").Value = "10:00-12:55 2:00-3:00"
$ws.Range("C2").Copy($ws.Range("C8"))
$ws.Range("C8").Value = 0.16319444444444445
$ws.Range("D8").Value = "Building the site"

# --- Add row 9: Dec 7, 2019 --------------------------------------------
$ws.Range("A2").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 43806
$ws.Range("B9").Value = "9:45-12:15 12:45-1:45"
$ws.Range("C2").Copy($ws.Range("C9"))
$ws.Range("C9").Value = 0.14583333333333334
$ws.Range("D9").Value = "Building site, setting up emails, setting up web host, making business cards"

# --- View state: un-zoom, scroll back to A1, select B10 ---------------
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B10").Select() | Out-Null
